$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with the default (unstyled) format, used to strip the
# quote-prefix style that Excel applies when a value is entered with a
# leading apostrophe (our 'force text' trick below).
$defaultStyle = $ws.Range("A1").Style

$ws.Range('D2').Formula = "'26.101.99"
$ws.Range('D2').Style = $defaultStyle
$ws.Range('E2').Formula = "'  -0.48%  "
$ws.Range('E2').Style = $defaultStyle
$ws.Range('D3').Formula = "'1.656.93"
$ws.Range('D3').Style = $defaultStyle
$ws.Range('E3').Formula = "'  -0.26%  "
$ws.Range('E3').Style = $defaultStyle
$ws.Range('D5').Formula = "'218.55"
$ws.Range('D5').Style = $defaultStyle
$ws.Range('E5').Formula = "'  -0.15%  "
$ws.Range('E5').Style = $defaultStyle
$ws.Range('D6').Formula = "'0.5297"
$ws.Range('D6').Style = $defaultStyle
$ws.Range('E6').Formula = "'  +1.38%  "
$ws.Range('E6').Style = $defaultStyle
$ws.Range('E7').Formula = "'  -0.29%  "
$ws.Range('E7').Style = $defaultStyle
$ws.Range('D8').Formula = "'0.2616"
$ws.Range('D8').Style = $defaultStyle
$ws.Range('E8').Formula = "'  -2.06%  "
$ws.Range('E8').Style = $defaultStyle
$ws.Range('D9').Formula = "'0.06337"
$ws.Range('D9').Style = $defaultStyle
$ws.Range('E9').Formula = "'  -0.12%  "
$ws.Range('E9').Style = $defaultStyle
$ws.Range('D10').Formula = "'20.42"
$ws.Range('D10').Style = $defaultStyle
$ws.Range('E10').Formula = "'  -3.11%  "
$ws.Range('E10').Style = $defaultStyle
$ws.Range('D11').Formula = "'0.07768"
$ws.Range('D11').Style = $defaultStyle
$ws.Range('E11').Formula = "'  +0.59%  "
$ws.Range('E11').Style = $defaultStyle
$ws.Range('B12').Formula = "'Polkadot"
$ws.Range('B12').Style = $defaultStyle
$ws.Range('C12').Formula = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range('C12').Style = $defaultStyle
$ws.Range('D12').Formula = "'4.498"
$ws.Range('D12').Style = $defaultStyle
$ws.Range('E12').Formula = "'  +1.51%  "
$ws.Range('E12').Style = $defaultStyle
$ws.Range('B13').Formula = "'WrappedEther"
$ws.Range('B13').Style = $defaultStyle
$ws.Range('C13').Formula = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range('C13').Style = $defaultStyle
$ws.Range('D13').Formula = "'1.657.11"
$ws.Range('D13').Style = $defaultStyle
$ws.Range('E13').Formula = "'  -0.31%  "
$ws.Range('E13').Style = $defaultStyle
$ws.Range('D14').Formula = "'0.5476"
$ws.Range('D14').Style = $defaultStyle
$ws.Range('E14').Formula = "'  -0.03%  "
$ws.Range('E14').Style = $defaultStyle
$ws.Range('D15').Formula = "'0.0₅8161"
$ws.Range('D15').Style = $defaultStyle
$ws.Range('E15').Formula = "'  -0.66%  "
$ws.Range('E15').Style = $defaultStyle
$ws.Range('D16').Formula = "'65.23"
$ws.Range('D16').Style = $defaultStyle
$ws.Range('E16').Formula = "'  +0.31%  "
$ws.Range('E16').Style = $defaultStyle
$ws.Range('D17').Formula = "'26.123.21"
$ws.Range('D17').Style = $defaultStyle
$ws.Range('E17').Formula = "'  -0.48%  "
$ws.Range('E17').Style = $defaultStyle
$ws.Range('E18').Formula = "'  -0.36%  "
$ws.Range('E18').Style = $defaultStyle
$ws.Range('D19').Formula = "'4.552"
$ws.Range('D19').Style = $defaultStyle
$ws.Range('E19').Formula = "'  -2.23%  "
$ws.Range('E19').Style = $defaultStyle
$ws.Range('D20').Formula = "'193.41"
$ws.Range('D20').Style = $defaultStyle
$ws.Range('E20').Formula = "'  -0.92%  "
$ws.Range('E20').Style = $defaultStyle
$ws.Range('D21').Formula = "'10.08"
$ws.Range('D21').Style = $defaultStyle
$ws.Range('E21').Formula = "'  -0.71%  "
$ws.Range('E21').Style = $defaultStyle
$ws.Range('D22').Formula = "'6.026"
$ws.Range('D22').Style = $defaultStyle
$ws.Range('E22').Formula = "'  -1.08%  "
$ws.Range('E22').Style = $defaultStyle
$ws.Range('E23').Formula = "'  -0.41%  "
$ws.Range('E23').Style = $defaultStyle
$ws.Range('D24').Formula = "'140.21"
$ws.Range('D24').Style = $defaultStyle
$ws.Range('E24').Formula = "'  +0.70%  "
$ws.Range('E24').Style = $defaultStyle
$ws.Range('D25').Formula = "'0.1244"
$ws.Range('D25').Style = $defaultStyle
$ws.Range('E25').Formula = "'  +0.09%  "
$ws.Range('E25').Style = $defaultStyle
$ws.Range('D26').Formula = "'7.274"
$ws.Range('D26').Style = $defaultStyle
$ws.Range('E26').Formula = "'  +0.51%  "
$ws.Range('E26').Style = $defaultStyle
$ws.Range('D27').Formula = "'16.16"
$ws.Range('D27').Style = $defaultStyle
$ws.Range('E27').Formula = "'  -0.34%  "
$ws.Range('E27').Style = $defaultStyle
$ws.Range('E28').Formula = "'  +1.12%  "
$ws.Range('E28').Style = $defaultStyle
$ws.Range('D29').Formula = "'0.05941"
$ws.Range('D29').Style = $defaultStyle
$ws.Range('E29').Formula = "'  -0.50%  "
$ws.Range('E29').Style = $defaultStyle
$ws.Range('D30').Formula = "'1.278"
$ws.Range('D30').Style = $defaultStyle
$ws.Range('E30').Formula = "'  -0.34%  "
$ws.Range('E30').Style = $defaultStyle
$ws.Range('D31').Formula = "'3.510"
$ws.Range('D31').Style = $defaultStyle
$ws.Range('E31').Formula = "'  -3.23%  "
$ws.Range('E31').Style = $defaultStyle
$ws.Range('D32').Formula = "'3.237"
$ws.Range('D32').Style = $defaultStyle
$ws.Range('E32').Formula = "'  -2.19%  "
$ws.Range('E32').Style = $defaultStyle
$ws.Range('D33').Formula = "'1.558"
$ws.Range('D33').Style = $defaultStyle
$ws.Range('E33').Formula = "'  -4.63%  "
$ws.Range('E33').Style = $defaultStyle
$ws.Range('E34').Formula = "'  -3.14%  "
$ws.Range('E34').Style = $defaultStyle
$ws.Range('E35').Formula = "'  -0.44%  "
$ws.Range('E35').Style = $defaultStyle
$ws.Range('D36').Formula = "'2.765"
$ws.Range('D36').Style = $defaultStyle
$ws.Range('E36').Formula = "'  -0.48%  "
$ws.Range('E36').Style = $defaultStyle
$ws.Range('D37').Formula = "'0.5640"
$ws.Range('D37').Style = $defaultStyle
$ws.Range('E37').Formula = "'  -4.43%  "
$ws.Range('E37').Style = $defaultStyle
$ws.Range('D38').Formula = "'0.01612"
$ws.Range('D38').Style = $defaultStyle
$ws.Range('E38').Formula = "'  +0.94%  "
$ws.Range('E38').Style = $defaultStyle
$ws.Range('D39').Formula = "'5.841"
$ws.Range('D39').Style = $defaultStyle
$ws.Range('E39').Formula = "'  -2.61%  "
$ws.Range('E39').Style = $defaultStyle
$ws.Range('D40').Formula = "'0.8491"
$ws.Range('D40').Style = $defaultStyle
$ws.Range('E40').Formula = "'  -0.99%  "
$ws.Range('E40').Style = $defaultStyle
$ws.Range('E41').Formula = "'  -0.27%  "
$ws.Range('E41').Style = $defaultStyle
$ws.Range('D42').Formula = "'101.46"
$ws.Range('D42').Style = $defaultStyle
$ws.Range('E42').Formula = "'  +1.54%  "
$ws.Range('E42').Style = $defaultStyle
$ws.Range('D43').Formula = "'1.012.10"
$ws.Range('D43').Style = $defaultStyle
$ws.Range('E43').Formula = "'  -1.72%  "
$ws.Range('E43').Style = $defaultStyle
$ws.Range('D44').Formula = "'1.801.94"
$ws.Range('D44').Style = $defaultStyle
$ws.Range('E44').Formula = "'  -0.06%  "
$ws.Range('E44').Style = $defaultStyle
$ws.Range('D45').Formula = "'57.10"
$ws.Range('D45').Style = $defaultStyle
$ws.Range('E45').Formula = "'  -0.49%  "
$ws.Range('E45').Style = $defaultStyle
$ws.Range('D46').Formula = "'1.002"
$ws.Range('D46').Style = $defaultStyle
$ws.Range('E46').Formula = "'  -0.40%  "
$ws.Range('E46').Style = $defaultStyle
$ws.Range('D47').Formula = "'0.0₈103"
$ws.Range('D47').Style = $defaultStyle
$ws.Range('E47').Formula = "'  -4.23%  "
$ws.Range('E47').Style = $defaultStyle
$ws.Range('E48').Formula = "'  +1.33%  "
$ws.Range('E48').Style = $defaultStyle
$ws.Range('D49').Formula = "'0.05156"
$ws.Range('D49').Style = $defaultStyle
$ws.Range('E49').Formula = "'  -0.61%  "
$ws.Range('E49').Style = $defaultStyle
$ws.Range('D50').Formula = "'1.467"
$ws.Range('D50').Style = $defaultStyle
$ws.Range('E50').Formula = "'  -0.05%  "
$ws.Range('E50').Style = $defaultStyle
$ws.Range('D51').Formula = "'7.733"
$ws.Range('D51').Style = $defaultStyle
$ws.Range('E51').Formula = "'  -4.28%  "
$ws.Range('E51').Style = $defaultStyle
